# Update column G ("K") values for rows 2-16 in the active sheet,
# per the diff: regen save_data to use K instead of Strike#,
# regen std/mean, calc and write s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 4
    6  = 0
    7  = 2
    8  = 1
    9  = 1
    10 = 1
    11 = 0
    12 = 1
    13 = 2
    14 = 1
    15 = 3
    16 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
